# Commit 9 June 2016
# Add a Standard Error (SE) row beneath the existing Mean/SD/Min/Max/N
# summary statistics block (columns K:N, row 7).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Label for the new stat row
$ws.Range("K7").Value = "SE"

# SE = SD / SQRT(N)  -- L3/M3/N3 hold SD, L6/M6/N6 hold N
$ws.Range("L7").Formula = "=L3/SQRT(L6)"
$ws.Range("M7:N7").Formula = "=M3/SQRT(M6)"

# Match the formatting used by the rest of the stats block (L2:N6)
$ws.Range("L7:N7").NumberFormat = $ws.Range("L2").NumberFormat
$ws.Range("L7:N7").HorizontalAlignment = $ws.Range("L2").HorizontalAlignment

# Update the active selection like the author left it
$ws.Range("M7").Select() | Out-Null
